$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1660.1428
$ws.Range("J17").Value = 1670.5186
$ws.Range("L17").Value = 5011.5558
$ws.Range("N17").Value = -5347.5558

$ws.Range("H18").Value = 3386
$ws.Range("J18").Value = 4160
$ws.Range("L18").Value = 4160
$ws.Range("N18").Value = -4728

$ws.Range("H31").Value = 5699.5
$ws.Range("I31").Value = 400
$ws.Range("K31").Value = 1200
$ws.Range("M31").Value = -970

$ws.Range("H32").Value = 7744
$ws.Range("J32").Value = 7744
$ws.Range("L32").Value = 7744
$ws.Range("N32").Value = -8396

$ws.Range("H38").Value = 2761
$ws.Range("I38").Value = 391.75
$ws.Range("K38").Value = 1175.25
$ws.Range("M38").Value = -803.25

$ws.Range("H92").Value = 3032.2173
$ws.Range("I92").Value = 1225.5883
$ws.Range("J92").Value = 8151
$ws.Range("K92").Value = 1225.5883
$ws.Range("L92").Value = 8151
$ws.Range("M92").Value = 22.41170000000011
$ws.Range("N92").Value = -10647

$ws.Range("H138").Value = 3118.634
$ws.Range("I138").Value = 1705.3572
$ws.Range("K138").Value = 5116.071599999999
$ws.Range("M138").Value = 23.92840000000069

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 685.4286
$ws.Range("I5").Value = 685.4286
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 685.4286
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = -573.4286
$ws.Range("M5").ClearContents()

$ws.Range("H32").Value = 6875.17
$ws.Range("I32").Value = 6213.7
$ws.Range("K32").Value = 6213.7
$ws.Range("M32").Value = -5926.7

$ws.Range("H45").Value = 7499.5
$ws.Range("I45").Value = 4999.5
$ws.Range("K45").Value = 4999.5
$ws.Range("M45").Value = -4622.5

$ws.Range("H102").Value = 8058.9287
$ws.Range("I102").Value = 6875.625
$ws.Range("K102").Value = 6875.625
$ws.Range("M102").Value = -5253.625

$ws.Range("H110").Value = 3704.5483
$ws.Range("I110").Value = 3316.182
$ws.Range("K110").Value = 3316.182
$ws.Range("M110").Value = -1271.182

$ws.Range("H123").Value = 172666.67
$ws.Range("J123").Value = 172666.67
$ws.Range("L123").Value = 172666.67
$ws.Range("N123").Value = -182466.67

$ws.Range("H132").Value = 2859269
$ws.Range("I132").Value = 2030.3667
$ws.Range("K132").Value = 6091.1001
$ws.Range("M132").Value = -3561.1001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 685.4286
$ws.Range("I4").Value = 685.4286
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 685.4286
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = -570.4286
$ws.Range("M4").ClearContents()

$ws.Range("H22").Value = 1193
$ws.Range("I22").Value = 1193
$ws.Range("K22").Value = 1193
$ws.Range("M22").Value = -1020

$ws.Range("H86").Value = 1591.3846
$ws.Range("I86").Value = 1345.2273
$ws.Range("J86").Value = 2945.25
$ws.Range("K86").Value = 1345.2273
$ws.Range("L86").Value = 2945.25
$ws.Range("M86").Value = -222.2273
$ws.Range("N86").Value = -5191.25

$ws.Range("H89").Value = 1591.3846
$ws.Range("I89").Value = 1345.2273
$ws.Range("J89").Value = 2945.25
$ws.Range("K89").Value = 6726.136500000001
$ws.Range("L89").Value = 14726.25
$ws.Range("M89").Value = -1110.136500000001
$ws.Range("N89").Value = -25958.25

$ws.Range("H99").Value = 3137.5
$ws.Range("I99").Value = 2366.6667
$ws.Range("K99").Value = 2366.6667
$ws.Range("M99").Value = -868.6667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43512548
$ws.Range("I31").Value = 62547572
$ws.Range("J31").Value = 3914.2856
$ws.Range("K31").Value = 62547572
$ws.Range("L31").Value = 3914.2856
$ws.Range("M31").Value = -62547277
$ws.Range("N31").Value = -4504.2856

$ws.Range("H34").Value = 43512548
$ws.Range("I34").Value = 62547572
$ws.Range("J34").Value = 3914.2856
$ws.Range("K34").Value = 62547572
$ws.Range("L34").Value = 3914.2856
$ws.Range("M34").Value = -62547370
$ws.Range("N34").Value = -4318.2856

$ws.Range("H58").Value = 3859.4348
$ws.Range("I58").Value = 3422.4375
$ws.Range("K58").Value = 3422.4375
$ws.Range("M58").Value = -3219.4375

$ws.Range("H136").Value = 3859.4348
$ws.Range("I136").Value = 3422.4375
$ws.Range("K136").Value = 10267.3125
$ws.Range("M136").Value = -7717.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 15015
$ws.Range("I3").Value = 8909
$ws.Range("K3").Value = 26727
$ws.Range("M3").Value = -26615

$ws.Range("H114").Value = 9590.75
$ws.Range("J114").Value = 12777.667
$ws.Range("L114").Value = 38333.001
$ws.Range("N114").Value = -44841.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 154.86667
$ws.Range("I2").Value = 150.75
$ws.Range("J2").Value = 171.33333
$ws.Range("K2").Value = 150.75
$ws.Range("L2").Value = 171.33333
$ws.Range("M2").Value = -37.75
$ws.Range("N2").Value = -397.33333

$ws.Range("H13").Value = 4393.625
$ws.Range("I13").Value = 625
$ws.Range("J13").Value = 5649.8335
$ws.Range("K13").Value = 625
$ws.Range("L13").Value = 5649.8335
$ws.Range("M13").Value = -486
$ws.Range("N13").Value = -5927.8335

$ws.Range("H70").Value = 10528
$ws.Range("I70").Value = 8312.571
$ws.Range("J70").Value = 12251.111
$ws.Range("K70").Value = 8312.571
$ws.Range("L70").Value = 12251.111
$ws.Range("M70").Value = -8042.571
$ws.Range("N70").Value = -12791.111

$ws.Range("H73").Value = 10528
$ws.Range("I73").Value = 8312.571
$ws.Range("J73").Value = 12251.111
$ws.Range("K73").Value = 8312.571
$ws.Range("L73").Value = 12251.111
$ws.Range("M73").Value = -7376.571
$ws.Range("N73").Value = -14123.111

$ws.Range("H97").Value = 6262.5
$ws.Range("I97").Value = 1509.5
$ws.Range("K97").Value = 1509.5
$ws.Range("M97").Value = -1013.5

$ws.Range("H102").Value = 2378.3
$ws.Range("I102").Value = 2475.889
$ws.Range("K102").Value = 2475.889
$ws.Range("M102").Value = -853.8890000000001

$ws.Range("H122").Value = 8939.529
$ws.Range("I122").Value = 7687.875
$ws.Range("J122").Value = 10052.111
$ws.Range("K122").Value = 23063.625
$ws.Range("L122").Value = 30156.333
$ws.Range("M122").Value = -20613.625
$ws.Range("N122").Value = -35056.333

$ws.Range("H138").Value = 99998.5
$ws.Range("J138").Value = 99998.5
$ws.Range("L138").Value = 99998.5
$ws.Range("N138").Value = -110278.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1726.1538
$ws.Range("I16").Value = 1634.75
$ws.Range("J16").Value = 2030.8334
$ws.Range("K16").Value = 1634.75
$ws.Range("L16").Value = 2030.8334
$ws.Range("M16").Value = -1464.75
$ws.Range("N16").Value = -2370.8334

$ws.Range("H22").Value = 7907.524
$ws.Range("I22").Value = 14849.5
$ws.Range("J22").Value = 3635.5386
$ws.Range("K22").Value = 14849.5
$ws.Range("L22").Value = 3635.5386
$ws.Range("M22").Value = -14554.5
$ws.Range("N22").Value = -4225.5386

$ws.Range("H27").Value = 7907.524
$ws.Range("I27").Value = 14849.5
$ws.Range("J27").Value = 3635.5386
$ws.Range("K27").Value = 14849.5
$ws.Range("L27").Value = 3635.5386
$ws.Range("M27").Value = -14742.5
$ws.Range("N27").Value = -3849.5386

$ws.Range("H55").Value = 1174.3334
$ws.Range("I55").Value = 459.7143
$ws.Range("J55").Value = 1799.625
$ws.Range("K55").Value = 459.7143
$ws.Range("L55").Value = 1799.625
$ws.Range("M55").Value = -286.7143
$ws.Range("N55").Value = -2145.625

$ws.Range("H58").Value = 49988.332
$ws.Range("I58").Value = 49983
$ws.Range("J58").Value = 49999
$ws.Range("K58").Value = 49983
$ws.Range("L58").Value = 49999
$ws.Range("N58").Value = -50519
$ws.Range("M58").Value = -49723

$ws.Range("H68").Value = 3791747.2
$ws.Range("J68").Value = 5976.8335
$ws.Range("L68").Value = 5976.8335
$ws.Range("N68").Value = -7474.8335

$ws.Range("H71").Value = 3791747.2
$ws.Range("J71").Value = 5976.8335
$ws.Range("L71").Value = 29884.1675
$ws.Range("N71").Value = -37372.1675

$ws.Range("H93").Value = 4277265
$ws.Range("I93").Value = 3849.3333
$ws.Range("J93").Value = 7940193
$ws.Range("K93").Value = 3849.3333
$ws.Range("L93").Value = 7940193
$ws.Range("M93").Value = -2601.3333
$ws.Range("N93").Value = -7942689

$ws.Range("H136").Value = 2126.4443
$ws.Range("I136").Value = 2126.4443
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6379.3329
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -3829.3329
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5305.9375
$ws.Range("I107").Value = 3761.3845
$ws.Range("K107").Value = 11284.1535
$ws.Range("M107").Value = -9364.1535

$ws.Range("H113").Value = 468.1154
$ws.Range("J113").Value = 698
$ws.Range("L113").Value = 2094
$ws.Range("N113").Value = -6434

$ws.Range("H132").Value = 348050.22
$ws.Range("I132").Value = 3329.4167
$ws.Range("K132").Value = 9988.250100000001
$ws.Range("M132").Value = -7458.250100000001

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").Value = 0
$ws.Range("L140").ClearContents()

$ws.Range("H141").Value = 84827.375
$ws.Range("J141").Value = 84827.375
$ws.Range("L141").Value = 84827.375
$ws.Range("N141").Value = -95187.375
